$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph and insert a
# new ListBullet paragraph right after it, containing the two instructors
# separated by a manual line break (matching the XML diff exactly: two
# runs, the first ending in <w:br/>, the second holding the trailing name).

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $candidate
        break
    }
}

# Find the paragraph's index within the Paragraphs collection so we can
# reference the freshly created siblings afterward.
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $target.Range.Start) {
        $idx = $i
        break
    }
}

# Create a new paragraph right after the heading.
$target.Range.InsertParagraphAfter()

# First new paragraph: first instructor line.
$p1 = $d.Paragraphs.Item($idx + 1)
$p1.Style = "ListBullet"
$p1.Range.Text = "5983729 - Fernando Vernilli Junior"

# Split off a second paragraph for the next instructor (keeps the two
# pieces of text in separate runs once merged back below).
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($idx + 2)
$p2.Style = "ListBullet"
$p2.Range.Text = "1922320 - Sebastiao Ribeiro"

# Insert a manual line break just before the paragraph mark that ends the
# first new paragraph (so it renders as <w:br/> at the end of its run).
$p1b = $d.Paragraphs.Item($idx + 1)
$r1b = $p1b.Range
$breakPos = $r1b.End - 1
$breakPoint = $d.Range($breakPos, $breakPos)
$breakPoint.Text = [char]11

# Remove the paragraph mark separating the two new paragraphs so the
# second instructor's text becomes a second run within the same
# paragraph (rather than its own paragraph).
$p1c = $d.Paragraphs.Item($idx + 1)
$r1c = $p1c.Range
$markStart = $r1c.End - 1
$markRange = $d.Range($markStart, $r1c.End)
$markRange.Delete()

Write-Output "Inserted instructors paragraph after index $idx"
